# Insert a new weekly data row for "Espinaca" (Mercado Mayorista Lo Valledor
# de Santiago) at row 594, pushing the existing rows 594-659 down to 595-660.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 594 - Excel shifts rows 594..659
# down to 595..660 and copies formatting from the row above (so column D
# keeps its date number format).
$ws.Rows.Item(594).Insert()

# Populate the newly inserted row 594 with the new record's values.
$ws.Cells.Item(594, 1).Value2 = 6
$ws.Cells.Item(594, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(594, 3).Value2 = "Metropolitana"
$ws.Cells.Item(594, 4).Value2 = 44918
$ws.Cells.Item(594, 5).Value2 = 13
$ws.Cells.Item(594, 6).Value2 = 100112012
$ws.Cells.Item(594, 7).Value2 = "Espinaca"
$ws.Cells.Item(594, 8).Value2 = "Sin especificar"
$ws.Cells.Item(594, 9).Value2 = "Primera"
$ws.Cells.Item(594, 10).Value2 = 630
$ws.Cells.Item(594, 11).Value2 = 6000
$ws.Cells.Item(594, 12).Value2 = 6500
$ws.Cells.Item(594, 13).Value2 = 6214
$ws.Cells.Item(594, 14).Value2 = "`$/cuna 10 kilos"
$ws.Cells.Item(594, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(594, 16).Value2 = 621
$ws.Cells.Item(594, 17).Value2 = 10
$ws.Cells.Item(594, 18).Value2 = "Hortaliza"
